$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Paragraphs.Item(1).Range.Text = "2026-02-14 Saturday"

$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "92÷5=18, 2"
$t.Cell(1, 2).Range.Text = "25÷4=6, 1"
$t.Cell(1, 3).Range.Text = "38÷9=4, 2"
$t.Cell(1, 4).Range.Text = "41÷6=6, 5"
$t.Cell(1, 5).Range.Text = "84÷8=10, 4"
$t.Cell(5, 1).Range.Text = "92÷3=30, 2"
$t.Cell(5, 2).Range.Text = "64÷9=7, 1"
$t.Cell(5, 3).Range.Text = "44÷3=14, 2"
$t.Cell(5, 4).Range.Text = "24÷6=4, 0"
$t.Cell(5, 5).Range.Text = "32÷2=16, 0"
$t.Cell(9, 1).Range.Text = "35÷6=5, 5"
$t.Cell(9, 2).Range.Text = "39÷5=7, 4"
$t.Cell(9, 3).Range.Text = "45÷2=22, 1"
$t.Cell(9, 4).Range.Text = "62÷8=7, 6"
$t.Cell(9, 5).Range.Text = "10÷9=1, 1"
$t.Cell(13, 1).Range.Text = "46÷3=15, 1"
$t.Cell(13, 2).Range.Text = "95÷8=11, 7"
$t.Cell(13, 3).Range.Text = "78÷5=15, 3"
$t.Cell(13, 4).Range.Text = "38÷9=4, 2"
$t.Cell(13, 5).Range.Text = "17÷7=2, 3"
$t.Cell(17, 1).Range.Text = "37÷8=4, 5"
$t.Cell(17, 2).Range.Text = "24÷6=4, 0"
$t.Cell(17, 3).Range.Text = "51÷9=5, 6"
$t.Cell(17, 4).Range.Text = "68÷6=11, 2"
$t.Cell(17, 5).Range.Text = "55÷6=9, 1"
